# Scheduled runner update: refresh computed leve-profit figures (market
# board price pulls) across the per-class profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 13335348
$ws.Range("I100").Value = 20834434
$ws.Range("J100").Value = 3641.4443
$ws.Range("K100").Value = 20834434
$ws.Range("L100").Value = 3641.4443
$ws.Range("M100").Value = -20833893
$ws.Range("N100").Value = -4723.4443

$ws.Range("H112").Value = 3371.0356
$ws.Range("I112").Value = 740
$ws.Range("J112").Value = 3943
$ws.Range("K112").Value = 2220
$ws.Range("L112").Value = 11829
$ws.Range("M112").Value = -1112
$ws.Range("N112").Value = -14045

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 394558.44
$ws.Range("I32").Value = 3047.028
$ws.Range("K32").Value = 3047.028
$ws.Range("M32").Value = -2760.028

$ws.Range("H55").Value = 29628.143
$ws.Range("J55").Value = 29628.143
$ws.Range("L55").Value = 29628.143
$ws.Range("N55").Value = -30258.143

$ws.Range("H80").Value = 11631
$ws.Range("I80").Value = 2250
$ws.Range("J80").Value = 16321.5
$ws.Range("K80").Value = 2250
$ws.Range("L80").Value = 16321.5
$ws.Range("M80").Value = -1252
$ws.Range("N80").Value = -18317.5

$ws.Range("H83").Value = 11631
$ws.Range("I83").Value = 2250
$ws.Range("J83").Value = 16321.5
$ws.Range("K83").Value = 6750
$ws.Range("L83").Value = 48964.5
$ws.Range("M83").Value = -1758
$ws.Range("N83").Value = -58948.5

$ws.Range("H88").Value = 3514.45
$ws.Range("I88").Value = 2100.2222
$ws.Range("J88").Value = 4671.5454
$ws.Range("K88").Value = 2100.2222
$ws.Range("L88").Value = 4671.5454
$ws.Range("M88").Value = -1694.2222
$ws.Range("N88").Value = -5483.5454

$ws.Range("H91").Value = 3514.45
$ws.Range("I91").Value = 2100.2222
$ws.Range("J91").Value = 4671.5454
$ws.Range("K91").Value = 2100.2222
$ws.Range("L91").Value = 4671.5454
$ws.Range("M91").Value = -696.2222000000002
$ws.Range("N91").Value = -7479.5454

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 13616.77
$ws.Range("I82").Value = 3449.8333
$ws.Range("J82").Value = 22331.285
$ws.Range("K82").Value = 3449.8333
$ws.Range("L82").Value = 22331.285
$ws.Range("M82").Value = -3066.8333
$ws.Range("N82").Value = -23097.285

$ws.Range("H85").Value = 13616.77
$ws.Range("I85").Value = 3449.8333
$ws.Range("J85").Value = 22331.285
$ws.Range("K85").Value = 3449.8333
$ws.Range("L85").Value = 22331.285
$ws.Range("M85").Value = -2123.8333
$ws.Range("N85").Value = -24983.285

$ws.Range("H86").Value = 2691.0303
$ws.Range("I86").Value = 2468.0454
$ws.Range("J86").Value = 3137
$ws.Range("K86").Value = 2468.0454
$ws.Range("L86").Value = 3137
$ws.Range("M86").Value = -1345.0454
$ws.Range("N86").Value = -5383

$ws.Range("H89").Value = 2691.0303
$ws.Range("I89").Value = 2468.0454
$ws.Range("J89").Value = 3137
$ws.Range("K89").Value = 12340.227
$ws.Range("L89").Value = 15685
$ws.Range("M89").Value = -6724.226999999999
$ws.Range("N89").Value = -26917

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3222.7058
$ws.Range("I31").Value = 1553.7894
$ws.Range("J31").Value = 8101.077
$ws.Range("K31").Value = 1553.7894
$ws.Range("L31").Value = 8101.077
$ws.Range("M31").Value = -1258.7894
$ws.Range("N31").Value = -8691.077000000001

$ws.Range("H34").Value = 3222.7058
$ws.Range("I34").Value = 1553.7894
$ws.Range("J34").Value = 8101.077
$ws.Range("K34").Value = 1553.7894
$ws.Range("L34").Value = 8101.077
$ws.Range("M34").Value = -1351.7894
$ws.Range("N34").Value = -8505.077000000001

$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").Value = $null

$ws.Range("H62").Value = 2937.7334
$ws.Range("I62").Value = 2960
$ws.Range("J62").Value = 2912.2856
$ws.Range("K62").Value = 2960
$ws.Range("L62").Value = 2912.2856
$ws.Range("M62").Value = -2336
$ws.Range("N62").Value = -4160.2856

$ws.Range("H65").Value = 2937.7334
$ws.Range("I65").Value = 2960
$ws.Range("J65").Value = 2912.2856
$ws.Range("K65").Value = 14800
$ws.Range("L65").Value = 14561.428
$ws.Range("M65").Value = -11680
$ws.Range("N65").Value = -20801.428

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 15760950
$ws.Range("I56").Value = 15760950
$ws.Range("K56").Value = 15760950
$ws.Range("M56").Value = -15760420

$ws.Range("H113").Value = 1000.6
$ws.Range("I113").Value = 844.25
$ws.Range("J113").Value = 1007.8721
$ws.Range("K113").Value = 2532.75
$ws.Range("L113").Value = 3023.6163
$ws.Range("M113").Value = -362.75
$ws.Range("N113").Value = -7363.6163

$ws.Range("H123").Value = 2499.6667
$ws.Range("I123").Value = 2500
$ws.Range("J123").Value = 2499.5
$ws.Range("K123").Value = 7500
$ws.Range("L123").Value = 7498.5
$ws.Range("M123").Value = -5050
$ws.Range("N123").Value = -12398.5

$ws.Range("H129").Value = 20834976
$ws.Range("I129").Value = 1925.7142
$ws.Range("J129").Value = 37038460
$ws.Range("K129").Value = 5777.142599999999
$ws.Range("L129").Value = 111115380
$ws.Range("M129").Value = -777.1425999999992
$ws.Range("N129").Value = -111125380

$ws.Range("H130").Value = 2000
$ws.Range("I130").Value = 2000
$ws.Range("K130").Value = 6000
$ws.Range("M130").Value = -980

$ws.Range("H131").Value = 815.48
$ws.Range("J131").Value = 844.1087
$ws.Range("L131").Value = 2532.3261
$ws.Range("N131").Value = -12612.3261

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 7106
$ws.Range("I15").Value = 7106
$ws.Range("K15").Value = 7106
$ws.Range("M15").Value = -6818

$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").Value = $null

$ws.Range("H41").Value = 9870.857
$ws.Range("I41").Value = 4342
$ws.Range("J41").Value = 10792.333
$ws.Range("K41").Value = 4342
$ws.Range("L41").Value = 10792.333
$ws.Range("M41").Value = -3952
$ws.Range("N41").Value = -11572.333

$ws.Range("H113").Value = 24390704
$ws.Range("I113").Value = 441.03226
$ws.Range("K113").Value = 1323.09678
$ws.Range("M113").Value = 846.9032199999999

